# Workbook edit: remove the "telepulesmax(telep)" column-block (old A:B),
# shifting the remaining three blocks left, then update a few cell texts
# in the (now) "max" block and the "bennevan(sor)" block, and finally move
# the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the first block of columns (A:B) entirely, shifting everything
#    else two columns to the left.
$ws.Range("A1:B1").EntireColumn.Delete()

# 2) Text edits inside the (now) "max" block, which used to be columns C:D
#    and is now columns A:B.
$ws.Range("B6").Value = "n * m"
$ws.Range("B7").Value = "ho[(i - 1) div m + 1, (i - 1) mod m + 1]"

# 3) Text edits inside the (now) "bennevan(sor)" block, which used to be
#    columns E:F and is now columns C:D.
$ws.Range("D3").Value = "j"
$ws.Range("D6").Value = "ho[sor, j] = maxho"

# 4) Move the active selection (cosmetic, matches the saved view state).
$ws.Range("J18").Select()
